# Auto-generated: apply scheduled market-price refresh updates to the Ifrit_Profits workbook.
# Each sheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) and the
# columns H-N hold the current market price / profit figures that the scheduled runner refreshes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1842.0667
$ws.Range("I62").Value = 1758.35
$ws.Range("J62").Value = 2009.5
$ws.Range("K62").Value = 1758.35
$ws.Range("L62").Value = 2009.5
$ws.Range("M62").Value = -1134.35
$ws.Range("N62").Value = -3257.5

$ws.Range("H65").Value = 1842.0667
$ws.Range("I65").Value = 1758.35
$ws.Range("J65").Value = 2009.5
$ws.Range("K65").Value = 8791.75
$ws.Range("L65").Value = 10047.5
$ws.Range("M65").Value = -5671.75
$ws.Range("N65").Value = -16287.5

$ws.Range("H97").Value = 2247.5
$ws.Range("J97").Value = 2247.5
$ws.Range("L97").Value = 6742.5
$ws.Range("N97").Value = -7734.5

$ws.Range("H112").Value = 55556980
$ws.Range("J112").Value = 55556980
$ws.Range("L112").Value = 166670940
$ws.Range("N112").Value = -166673156

$ws.Range("H129").Value = 1389.421
$ws.Range("J129").Value = 1511.0625
$ws.Range("L129").Value = 4533.1875
$ws.Range("N129").Value = -14533.1875

$ws.Range("H137").Value = 1649.9
$ws.Range("I137").Value = 1143.9688
$ws.Range("J137").Value = 2228.1072
$ws.Range("K137").Value = 3431.9064
$ws.Range("L137").Value = 6684.321599999999
$ws.Range("M137").Value = -881.9064000000003
$ws.Range("N137").Value = -11784.3216

$ws.Range("H138").Value = 969108.2
$ws.Range("I138").Value = 2502.0715
$ws.Range("J138").Value = 1557477.1
$ws.Range("K138").Value = 7506.2145
$ws.Range("L138").Value = 4672431.300000001
$ws.Range("M138").Value = -2366.2145
$ws.Range("N138").Value = -4682711.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8276.392
$ws.Range("I32").Value = 3501.5662
$ws.Range("K32").Value = 3501.5662
$ws.Range("M32").Value = -3214.5662

$ws.Range("H58").Value = 28000
$ws.Range("J58").Value = 28000
$ws.Range("L58").Value = 28000
$ws.Range("N58").Value = -28860

$ws.Range("H74").Value = 3641.1765
$ws.Range("I74").Value = 641.6799999999999
$ws.Range("J74").Value = 11973.111
$ws.Range("K74").Value = 641.6799999999999
$ws.Range("L74").Value = 11973.111
$ws.Range("M74").Value = 232.3200000000001
$ws.Range("N74").Value = -13721.111

$ws.Range("H77").Value = 3641.1765
$ws.Range("I77").Value = 641.6799999999999
$ws.Range("J77").Value = 11973.111
$ws.Range("K77").Value = 3208.4
$ws.Range("L77").Value = 59865.55500000001
$ws.Range("M77").Value = 1159.6
$ws.Range("N77").Value = -68601.55500000001

$ws.Range("H126").Value = 5216
$ws.Range("I126").Value = 5216
$ws.Range("K126").Value = 15648
$ws.Range("M126").Value = -13178

$ws.Range("H132").Value = 1151148.2
$ws.Range("I132").Value = 1278905.5
$ws.Range("K132").Value = 3836716.5
$ws.Range("M132").Value = -3834186.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1060.8077
$ws.Range("I20").Value = 798.86664
$ws.Range("K20").Value = 798.86664
$ws.Range("M20").Value = -551.86664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1411.509
$ws.Range("I31").Value = 1100
$ws.Range("J31").Value = 1813.875
$ws.Range("K31").Value = 1100
$ws.Range("L31").Value = 1813.875
$ws.Range("M31").Value = -805
$ws.Range("N31").Value = -2403.875

$ws.Range("H34").Value = 1411.509
$ws.Range("I34").Value = 1100
$ws.Range("J34").Value = 1813.875
$ws.Range("K34").Value = 1100
$ws.Range("L34").Value = 1813.875
$ws.Range("M34").Value = -898
$ws.Range("N34").Value = -2217.875

$ws.Range("H58").Value = 3175.861
$ws.Range("I58").Value = 1639.7858
$ws.Range("J58").Value = 4153.364
$ws.Range("K58").Value = 1639.7858
$ws.Range("L58").Value = 4153.364
$ws.Range("M58").Value = -1436.7858
$ws.Range("N58").Value = -4559.364

$ws.Range("H105").Value = 937.5
$ws.Range("I105").Value = 1020
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 1020
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 727
$ws.Range("N105").Value = -4294

$ws.Range("H136").Value = 3175.861
$ws.Range("I136").Value = 1639.7858
$ws.Range("J136").Value = 4153.364
$ws.Range("K136").Value = 4919.357400000001
$ws.Range("L136").Value = 12460.092
$ws.Range("M136").Value = -2369.357400000001
$ws.Range("N136").Value = -17560.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1308.1818
$ws.Range("I14").Value = 1308.1818
$ws.Range("K14").Value = 3924.5454
$ws.Range("M14").Value = -3751.5454

$ws.Range("H133").Value = 5828.357
$ws.Range("I133").Value = 2524.4443
$ws.Range("K133").Value = 7573.3329
$ws.Range("M133").Value = -2513.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2067.375
$ws.Range("I7").Value = 2067.375
$ws.Range("K7").Value = 2067.375
$ws.Range("M7").Value = -1955.375

$ws.Range("H16").Value = 1776.9565
$ws.Range("I16").Value = 2159.647
$ws.Range("J16").Value = 692.6667
$ws.Range("K16").Value = 2159.647
$ws.Range("L16").Value = 692.6667
$ws.Range("M16").Value = -1989.647
$ws.Range("N16").Value = -1032.6667

$ws.Range("H22").Value = 292.5
$ws.Range("I22").Value = 305.55554
$ws.Range("J22").Value = 253.33333
$ws.Range("K22").Value = 305.55554
$ws.Range("L22").Value = 253.33333
$ws.Range("M22").Value = -10.55554000000001
$ws.Range("N22").Value = -843.3333299999999

$ws.Range("H27").Value = 292.5
$ws.Range("I27").Value = 305.55554
$ws.Range("J27").Value = 253.33333
$ws.Range("K27").Value = 305.55554
$ws.Range("L27").Value = 253.33333
$ws.Range("M27").Value = -198.55554
$ws.Range("N27").Value = -467.33333

$ws.Range("H46").Value = 845.5454999999999
$ws.Range("I46").Value = 771.5714
$ws.Range("J46").Value = 975
$ws.Range("K46").Value = 771.5714
$ws.Range("L46").Value = 975
$ws.Range("M46").Value = -583.5714
$ws.Range("N46").Value = -1351

$ws.Range("H55").Value = 160.6842
$ws.Range("I55").Value = 151.1
$ws.Range("J55").Value = 171.33333
$ws.Range("K55").Value = 151.1
$ws.Range("L55").Value = 171.33333
$ws.Range("M55").Value = 21.90000000000001
$ws.Range("N55").Value = -517.3333299999999

$ws.Range("H61").Value = 2984
$ws.Range("I61").Value = 1580.8
$ws.Range("K61").Value = 1580.8
$ws.Range("M61").Value = -1378.8

$ws.Range("H82").Value = 2666.6667
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("N82").Value = -3722

$ws.Range("H85").Value = 2666.6667
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 3000
$ws.Range("N85").Value = -5496

$ws.Range("H93").Value = 2177.2415
$ws.Range("I93").Value = 1610.6471
$ws.Range("K93").Value = 1610.6471
$ws.Range("M93").Value = -362.6470999999999

$ws.Range("H100").Value = 2275.5
$ws.Range("I100").Value = 2180
$ws.Range("J100").Value = 2434.6667
$ws.Range("K100").Value = 2180
$ws.Range("L100").Value = 2434.6667
$ws.Range("M100").Value = -1639
$ws.Range("N100").Value = -3516.6667

$ws.Range("H113").Value = 2984
$ws.Range("I113").Value = 1580.8
$ws.Range("K113").Value = 1580.8
$ws.Range("M113").Value = 589.2

$ws.Range("H126").Value = 2067.375
$ws.Range("I126").Value = 2067.375
$ws.Range("K126").Value = 6202.125
$ws.Range("M126").Value = -3732.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1011.5833
$ws.Range("I136").Value = 946.619
$ws.Range("K136").Value = 2839.857
$ws.Range("M136").Value = -289.857

